$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("weibull")
$ws1.Range("B2").Value = -2.61351770303824
$ws1.Range("C2").Value = 0.0912646059721666
$ws1.Range("B3").Value = -0.0447916602245716
$ws1.Range("C3").Value = 0.0750187544228287

$ws2 = $wb.Worksheets.Item("lognormal")
$ws2.Range("B2").Value = 2.37960527088054
$ws2.Range("C2").Value = 0.179786460979122
$ws2.Range("B3").Value = -1.07901445539707
$ws2.Range("C3").Value = 0.10031066104669

$ws3 = $wb.Worksheets.Item("llogis")
$ws3.Range("B2").Value = -2.12776853505119
$ws3.Range("C2").Value = 0.0802505639146212
$ws3.Range("B3").Value = 0.682229582156072
$ws3.Range("C3").Value = 0.102850940330416

$ws4 = $wb.Worksheets.Item("gompertz")
$ws4.Range("B2").Value = -2.30353901599517
$ws4.Range("C2").Value = 0.0849990247310865
$ws4.Range("B3").Value = -0.0295289955002821
$ws4.Range("C3").Value = 0.00740228834920616

$ws6 = $wb.Worksheets.Item("weibull cov")
$ws6.Range("A2").Value = 0.00832922830325483
$ws6.Range("B2").Value = -0.00329069460677465
$ws6.Range("A3").Value = -0.00329069460677465
$ws6.Range("B3").Value = 0.00562781351515267

$ws7 = $wb.Worksheets.Item("lognormal cov")
$ws7.Range("A2").Value = 0.0323231715513972
$ws7.Range("B2").Value = -0.0162601111896603
$ws7.Range("A3").Value = -0.0162601111896603
$ws7.Range("B3").Value = 0.010062228719624

$ws8 = $wb.Worksheets.Item("llogis cov")
$ws8.Range("A2").Value = 0.0064401530086147
$ws8.Range("B2").Value = 0.00421831753931453
$ws8.Range("A3").Value = 0.00421831753931453
$ws8.Range("B3").Value = 0.0105783159268507

$ws9 = $wb.Worksheets.Item("gompertz cov")
$ws9.Range("A2").Value = 0.00722483420523586
$ws9.Range("B2").Value = -0.000207469315286281
$ws9.Range("A3").Value = -0.000207469315286281
$ws9.Range("B3").Value = 0.0000547938728047933
